$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1155
$ws1.Range("F5").Value = 1045
$ws1.Range("F6").Value = 1836
$ws1.Range("F7").Value = 583
$ws1.Range("F8").Value = 1219
$ws1.Range("F10").Value = 14
$ws1.Range("F12").Value = 312
$ws1.Range("F13").Value = 86
$ws1.Range("F15").Value = 728
$ws1.Range("F16").Value = 190
$ws1.Range("F20").Value = 333
$ws1.Range("F21").Value = 174
$ws1.Range("F22").Value = 684
$ws1.Range("F23").Value = 50
$ws1.Range("F24").Value = 654
$ws1.Range("F25").Value = 170
$ws1.Range("F26").Value = 40
$ws1.Range("F27").Value = 885
$ws1.Range("F28").Value = 326
$ws1.Range("F29").Value = 169
$ws1.Range("F30").Value = 52
$ws1.Range("F31").Value = 289
$ws1.Range("F33").Value = 19
$ws1.Range("F34").Value = 412

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 324
$ws2.Range("F7").Value = 262
$ws2.Range("F11").Value = 122
$ws2.Range("F12").Value = 23

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1344
$ws4.Range("F5").Value = 1155
$ws4.Range("F6").Value = 1046
$ws4.Range("F7").Value = 1836
$ws4.Range("F8").Value = 583
$ws4.Range("F9").Value = 1219
$ws4.Range("F12").Value = 14
$ws4.Range("F14").Value = 312
$ws4.Range("F15").Value = 86
$ws4.Range("F17").Value = 728
$ws4.Range("F18").Value = 190
$ws4.Range("F22").Value = 324
$ws4.Range("F25").Value = 333
$ws4.Range("F27").Value = 262
$ws4.Range("F28").Value = 262
$ws4.Range("F29").Value = 174
$ws4.Range("F30").Value = 684
$ws4.Range("F31").Value = 50
$ws4.Range("F32").Value = 654
$ws4.Range("F33").Value = 170
$ws4.Range("F34").Value = 40
$ws4.Range("F35").Value = 885
$ws4.Range("F36").Value = 326
$ws4.Range("F39").Value = 169
$ws4.Range("F40").Value = 52
$ws4.Range("F41").Value = 289
$ws4.Range("F43").Value = 122
$ws4.Range("F44").Value = 122
$ws4.Range("F46").Value = 19
$ws4.Range("F47").Value = 23
$ws4.Range("F48").Value = 412
